$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 with license-related dataset mapping entry
$ws.Range("B17").Value = "field_license_wbddh"
$ws.Range("A17").Value = "license_title"
$ws.Range("C17").Value = $true

# Update selection to mirror the authored change (active cell now A17)
$ws.Range("A17").Select()
